$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "plain" decimal number (e.g. "327.30") must be
# forced to Text so Excel does not silently coerce them to Number and drop
# the trailing zero (the source data keeps these as text in column D).
$textCells = @("D5","D7","D9","D10","D11","D13","D14","D15","D16","D22","D23","D24","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48")

$updates = [ordered]@{
    "D2" = "28.750.68"
    "E2" = "  +1.70%  "
    "D3" = "1.871.03"
    "E3" = "  +1.76%  "
    "E4" = "  +0.47%  "
    "D5" = "327.30"
    "E5" = "  -0.75%  "
    "E6" = "  +0.39%  "
    "D7" = "0.4645"
    "E7" = "  +1.10%  "
    "E8" = "  +1.49%  "
    "D9" = "0.07907"
    "E9" = "  +0.55%  "
    "D10" = "0.9711"
    "E10" = "  +0.98%  "
    "D11" = "22.32"
    "E11" = "  +1.68%  "
    "D12" = "1.891.67"
    "E12" = "  +2.20%  "
    "D13" = "5.745"
    "E13" = "  +0.46%  "
    "D14" = "6.937"
    "E14" = "  +0.33%  "
    "D15" = "0.06965"
    "E15" = "  +1.66%  "
    "D16" = "88.36"
    "E16" = "  +1.68%  "
    "E17" = "  +0.43%  "
    "E18" = "  +1.37%  "
    "E19" = "  +0.25%  "
    "D21" = "28.754.70"
    "E21" = "  +1.62%  "
    "D22" = "5.332"
    "D23" = "11.10"
    "E23" = "  +1.27%  "
    "D24" = "2.121"
    "E24" = "  -0.87%  "
    "D25" = "2.109.72"
    "E25" = "  +3.47%  "
    "D26" = "153.58"
    "E26" = "  +0.05%  "
    "D27" = "19.38"
    "E27" = "  +0.86%  "
    "D28" = "5.728"
    "E28" = "  -0.12%  "
    "D29" = "2.003"
    "E29" = "  +0.91%  "
    "D30" = "119.71"
    "E30" = "  +2.33%  "
    "D31" = "0.09373"
    "E31" = "  +0.72%  "
    "D32" = "0.9364"
    "E32" = "  -0.68%  "
    "D33" = "5.331"
    "E33" = "  +0.89%  "
    "D34" = "1.349"
    "E34" = "  +1.84%  "
    "D35" = "3.360"
    "D36" = "0.05849"
    "E36" = "  -2.77%  "
    "D37" = "0.02133"
    "E37" = "  -0.76%  "
    "D38" = "1.150"
    "E38" = "  +0.47%  "
    "D39" = "7.922"
    "E39" = "  +3.93%  "
    "D40" = "0.5667"
    "E40" = "  +0.71%  "
    "D41" = "9.980"
    "E41" = "  -0.18%  "
    "D42" = "0.1786"
    "E42" = "  +0.41%  "
    "D43" = "0.07243"
    "E43" = "  +3.06%  "
    "D44" = "11.76"
    "E44" = "  +0.74%  "
    "D45" = "0.5321"
    "E45" = "  +0.56%  "
    "B46" = "WEMIXToken"
    "C46" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D46" = "1.138"
    "E46" = "  -9.19%  "
    "B47" = "RenderToken"
    "C47" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D47" = "2.131"
    "E47" = "  -6.57%  "
    "D48" = "1.852"
    "E48" = "  +1.03%  "
    "E49" = "  +1.14%  "
    "E50" = "  +0.69%  "
    "E51" = "  +0.50%  "
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    if ($textCells -contains $cellRef) {
        $range.NumberFormat = "@"
        $range.Value = $updates[$cellRef]
        $range.ClearFormats()
    } else {
        $range.Value = $updates[$cellRef]
    }
}
